$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.22%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.27%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.104"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.91%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08049"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.37%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.942"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-13.45%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.996"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.85%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9327"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.68%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1454"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.22%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1924"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.08%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09016"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.28%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03502"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.54%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09789"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.22%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001394"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.49%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006044"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-4.36%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.778"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.61%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.197"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.20%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3420"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.92%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1328"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.29%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.789"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.19%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2413"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.04%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04394"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.77%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001239"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.69%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004271"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.69%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.14%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02044"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.84%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05028"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.38%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007437"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.69%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01013"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.12%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1349"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.91%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002122"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.26%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009074"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-8.68%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006199"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.27%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.23%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002803"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "28.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.23%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.23%"
